$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G"  = -0.3093617021276596
    "H"  = -0.3093617021276596
    "I"  = -0.2914893617021277
    "J"  = -0.2914893617021277
    "K"  = -0.738
    "L"  = -0.3140425531914893
    "U"  = 5.78
    "V"  = 0.4587301587301588
    "W"  = -0.4122905027932961
    "X"  = 0.02039890644709087
    "Y"  = -0.432689409240387
    "Z"  = 6.714285714285713
    "AA" = -1.957142857142857
    "AB" = 0.02113518044582894
    "AC" = -1.978278037588686
    "AD" = 2.21
    "AF" = 2.21
    "AG" = -3.57
    "AH" = 0.149223497636732
    "AI" = 0.6820987654320987
    "AJ" = -0.3953488372093024
    "AK" = 1.405511811023622
    "AL" = 0.005
    "AM" = -0.009000000000000001
    "AO" = -137
    "AQ" = 76.11111111111111
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
